$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    3 = @(36596, 5374, 6044)
    4 = @(18543, 2158, 2065)
    5 = @(62881, 5513, 5774)
    6 = @(1424, 514, 133)
    7 = @(39951, 6573, 5380)
    8 = @(4742, 1075, 972)
    9 = @(4918, 983, 590)
    10 = @(2332, 329, 228)
    11 = @(332, 168, 29)
    12 = @(3, 0, 0)
    13 = @(732, 193, 228)
    14 = @(2509, 1051, 738)
    15 = @(4297, 1596, 705)
    16 = @(2718, 1070, 406)
    17 = @(1787, 695, 137)
    18 = @(14485, 2292, 2690)
    19 = @(1330, 515, 401)
    20 = @(15796, 1877, 2634)
    21 = @(140, 302, 10)
    22 = @(14304, 1910, 2591)
    23 = @(921, 323, 130)
    24 = @(16192, 2350, 3187)
    25 = @(65503, 6254, 7928)
    26 = @(4754, 1541, 833)
    27 = @(0, 0, 0)
    28 = @(4523, 943, 1065)
    29 = @(1203, 389, 239)
    30 = @(12460, 2239, 2144)
    31 = @(360, 160, 181)
    32 = @(2279, 1431, 338)
    33 = @(14164, 2784, 2232)
    34 = @(8709, 2537, 1837)
    35 = @(5023, 529, 1258)
    36 = @(47314, 4749, 4678)
    37 = @(7067, 2342, 1039)
    38 = @(20635, 1684, 2414)
    39 = @(825, 782, 174)
    40 = @(1897, 404, 681)
    41 = @(2228, 271, 93)
    42 = @(8333, 482, 265)
    43 = @(236, 89, 66)
    44 = @(536, 38, 43)
    45 = @(1045, 14, 7)
    46 = @(2916, 764, 366)
    47 = @(9811, 2928, 1749)
    48 = @(27143, 2909, 3731)
    49 = @(12563, 3050, 1041)
    50 = @(9334, 923, 1309)
    51 = @(26840, 2495, 3834)
    52 = @(4056, 458, 1077)
    53 = @(12340, 2553, 2061)
    54 = @(1493, 1041, 601)
    55 = @(1739, 1073, 128)
    56 = @(3031, 733, 841)
    57 = @(10256, 3910, 2077)
    58 = @(12378, 894, 474)
    59 = @(567293, 85450, 79325)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
